$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-10 (columns D, K, L, M, N, O, P, Q, R, S, T)
# This represents a reordering/update of the weekly price records.
$data = @(
    @{ Row = 2;  D = 44545; K = "Castle Brite"; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; Q = "$/caja 15 kilos";        R = "Región de O'Higgins"; S = 1233; T = 15 },
    @{ Row = 3;  D = 44545; K = "Castle Brite"; L = "Segunda"; M = 50;  N = 17000; O = 17000; P = 17000; Q = "$/caja 15 kilos";        R = "Región de O'Higgins"; S = 1133; T = 15 },
    @{ Row = 4;  D = 44559; K = "Modesto";      L = "Primera"; M = 100; N = 19000; O = 20000; P = 19500; Q = "$/caja 18 kilos";        R = "Región de O'Higgins"; S = 1083; T = 18 },
    @{ Row = 5;  D = 44559; K = "Modesto";      L = "Segunda"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "$/caja 18 kilos";        R = "Región de O'Higgins"; S = 1000; T = 18 },
    @{ Row = 6;  D = 44187; K = "Dina";         L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "$/caja 18 kilos";        R = "Región Metropolitana"; S = 861;  T = 18 },
    @{ Row = 7;  D = 44159; K = "Castle Brite"; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "$/caja 15 kilos";        R = "Región Metropolitana"; S = 967;  T = 15 },
    @{ Row = 8;  D = 44579; K = "Modesto";      L = "Primera"; M = 180; N = 13000; O = 14000; P = 13444; Q = "$/caja 18 kilos";        R = "Región Metropolitana"; S = 747;  T = 18 },
    @{ Row = 9;  D = 44189; K = "Dina";         L = "Primera"; M = 200; N = 15000; O = 16000; P = 15500; Q = "$/caja 15 kilos granel"; R = "Región de O'Higgins"; S = 1033; T = 15 },
    @{ Row = 10; D = 44189; K = "Dina";         L = "Segunda"; M = 100; N = 14000; O = 14000; P = 14000; Q = "$/caja 15 kilos granel"; R = "Región de O'Higgins"; S = 933;  T = 15 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 4).Value  = $rec.D   # D - Fecha
    $ws.Cells.Item($r, 11).Value = $rec.K   # K - Variedad
    $ws.Cells.Item($r, 12).Value = $rec.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $rec.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $rec.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $rec.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $rec.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $rec.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $rec.R   # R - Origen
    $ws.Cells.Item($r, 19).Value = $rec.S   # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $rec.T   # T - Kg / unidad
}
